$d = $word.ActiveDocument
$r = $d.Range(0, 2)
$d.Bookmarks.Add("TEST_0_2", $r)
